# Update "想去人数" (interest count) figures for the latest scrape run.
# The same events are mirrored across the "展览" and "全部类型" sheets,
# but row numbers differ slightly between the two (展) sheet has one
# fewer row above row 23 than the "全部类型" sheet), so each sheet gets
# its own explicit row -> new-value mapping.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value  = 14146
$ws1.Range("F7").Value  = 15797
$ws1.Range("F8").Value  = 11
$ws1.Range("F9").Value  = 62
$ws1.Range("F21").Value = 129
$ws1.Range("F23").Value = 17
$ws1.Range("F27").Value = 5591
$ws1.Range("F28").Value = 77
$ws1.Range("F32").Value = 6

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value  = 14146
$ws4.Range("F7").Value  = 15797
$ws4.Range("F8").Value  = 11
$ws4.Range("F9").Value  = 62
$ws4.Range("F21").Value = 129
$ws4.Range("F24").Value = 17
$ws4.Range("F28").Value = 5591
$ws4.Range("F29").Value = 77
$ws4.Range("F33").Value = 6
